# edit.ps1 - applies the resume edits described by the target diff.
$d = $word.ActiveDocument
$lf = [char]11   # manual line break (<w:br/>) character, as seen in Range.Text

function Replace-Text($find, $replace, [bool]$matchCase = $true) {
    $d.Content.Find.Execute($find, $matchCase, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) "Highly motivated and believes in taking initiative"
#    -> "Highly motivated and believe in taking initiative"
# ---------------------------------------------------------------------------
Replace-Text "Highly motivated and believes in taking initiative" `
             "Highly motivated and believe in taking initiative"

# ---------------------------------------------------------------------------
# 2) Update the personal-site hyperlink address/display text, then add a
#    trailing ")" right after it.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Address -eq "https://bit.ly/vitalkravitz") {
        $h.TextToDisplay = "https://bit.ly/vitalkravitz-new"
        $h.Address = "https://bit.ly/vitalkravitz-new"
        $hEnd = $h.Range.End
        $ins = $d.Range($hEnd, $hEnd)
        $ins.InsertAfter(")")
        $ins.Bold = 0
        $ins.BoldBi = 0
        break
    }
}

# ---------------------------------------------------------------------------
# 3) "Manual QA testing and Automation" -> "Manual QA testing, and Automation"
# ---------------------------------------------------------------------------
Replace-Text "Manual QA testing and" "Manual QA testing, and"

# ---------------------------------------------------------------------------
# 4) "Very good understanding in QA methodologies" -> "...understanding of QA..."
# ---------------------------------------------------------------------------
Replace-Text "Very good understanding in QA methodologies" `
             "Very good understanding of QA methodologies"

# ---------------------------------------------------------------------------
# 5) After "CSS", add " (another site I've built <link>)" before the break
#    that starts the "Working with the next test management tools" line,
#    then drop "the next " from that line's text.
# ---------------------------------------------------------------------------
$cssRng = $d.Content
$cssRng.Find.Execute("CSS") | Out-Null
$cssRng.Collapse(0) | Out-Null
$insCss = $d.Range($cssRng.End, $cssRng.End)
$insCss.InsertAfter(" (another site I've built PLACEHOLDERZANZIBARLINK)")
$insCss.Bold = 0
$insCss.BoldBi = 0

$linkRng = $d.Content
$linkRng.Find.Execute("PLACEHOLDERZANZIBARLINK") | Out-Null
$d.Hyperlinks.Add($linkRng, "https://bit.ly/new-zanzibar", "", "", "https://bit.ly/new-zanzibar") | Out-Null

Replace-Text "Working with the next test management tools: " `
             "Working with test management tools: "

# ---------------------------------------------------------------------------
# 6) Drop the whole "Working with the next OS Platform: windows, Linux" line
#    (and its manual line break), then tweak the following line:
#      "Working with variety of testing tools: ..."
#        -> "Working with a variety of testing tools: ..."
#      "...DevTools and Firefox Developer Tools"
#        -> "...DevTool and Firefox Developer Tool"
# ---------------------------------------------------------------------------
$searchOs = "testing" + $lf + "Working with the next OS Platform: windows, Linux" + $lf + "Working with variety"
$replaceOs = "testing" + $lf + "Working with variety"
Replace-Text $searchOs $replaceOs

Replace-Text "Working with variety of testing tools: " `
             "Working with a variety of testing tools: "

Replace-Text "DevTools" "DevTool"
Replace-Text "Firefox Developer Tools" "Firefox Developer Tool"

# ---------------------------------------------------------------------------
# 7) "graphic design and digital marketing" -> "graphic design, and digital marketing"
# ---------------------------------------------------------------------------
Replace-Text "esign and" "esign, and"

# ---------------------------------------------------------------------------
# 8) "In charge of the of entire food, beverages and Consumable's"
#    -> "In charge of the entire food, beverages, and Consumable's"
# ---------------------------------------------------------------------------
Replace-Text "of entire food, beverages and" "entire food, beverages, and"

# ---------------------------------------------------------------------------
# 9) "With emphasis on FC demands" -> "With an emphasis on FC demands"
# ---------------------------------------------------------------------------
Replace-Text "ith emphasis on" "ith an emphasis on"

# ---------------------------------------------------------------------------
# 10) "according of all facility departments" -> "according to of all facility departments"
# ---------------------------------------------------------------------------
Replace-Text "according of all facility" "according to of all facility"

# ---------------------------------------------------------------------------
# 11) "Windows (all), Linux, MacOS, Mobile (iOS, Android)"
#     -> "Windows (all), Linux, macOS, Mobile (iOS, Android)"
# ---------------------------------------------------------------------------
Replace-Text "MacOS" "macOS"

# ---------------------------------------------------------------------------
# 12) "Fluent Speaking, Reading and Writing with basic competence"
#     -> "Fluent Speaking, Reading, and Writing with basic competence"
# ---------------------------------------------------------------------------
Replace-Text "Reading and Writing with basic competence" `
             "Reading, and Writing with basic competence"

Write-Output "edit.ps1 completed"
